$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "typ"
$ws.Range("B1").Value = "ks"
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = "B"
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = 5
